$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the spelling mistake in the branch name "Навоийскиц филиал" -> "Навоийский филиал"
$ws.Range("B11").Value = "Навоийский филиал"
